$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
Write-Output $nm
$t = $nm.Theme
Write-Output $t
$cs = $t.ThemeColorScheme
for ($i=1; $i -le 12; $i++) {
  $c = $cs.Item($i)
  $rgb = $c.RGB
  Write-Output "Item $i : $rgb"
}
